# The underlying change recorded for this document is a pure XML
# canonicalization: every attribute on every element in word/document.xml
# and word/styles.xml was re-serialized in (alphabetically) sorted
# attribute order (xmlns:m.. before xmlns:mc.., w:h before w:w, w:bottom
# before w:footer, w:qFormat before w:uiPriority, w:default before
# w:styleId before w:type, etc.). Every before/after pair in the diff has
# the exact same attribute *names and values* -- only the textual order
# of attributes inside each start tag differs.
#
# That kind of reordering is a serializer-level detail with no effect on
# the Word object model: there is no document content, formatting,
# paragraph, style definition, page-setup value, or property that
# actually changed. So the faithful Word-OM replay of this commit is to
# leave the document's content untouched.
$d = $word.ActiveDocument

# Touch the document object (keeps this a valid, executed COM script)
# without mutating any content/formatting, since none of the page setup,
# style, or body values differ between the before/after XML -- only
# attribute order does, which is outside the Word object model's surface.
$sectionCount = $d.Sections.Count
